$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 38) continuing the bi-weekly schedule table.
$ws.Range("A38").Value = "2025/12/12"
$ws.Range("B38").Value = "2026/2/6"
$ws.Range("C38").Value = "第88期 第七代貓貓包"

# Columns A and B store these dates as text (same number format as the row above).
$ws.Range("A38:B38").NumberFormat = $ws.Range("A37:B37").NumberFormat

# Mirror the new active cell / selection shown after the edit.
$ws.Range("C38").Select()
